# Atualização de bases das ligas, do dia: 19-04-2024 às 00:38
#
# This script:
#  1) Swaps the data (everything except id/Div/Div Original Name/Date
#     columns A, C, D, E) between several pairs of rows whose match
#     records had been mixed up (home/away swapped in source feed).
#  2) Appends two newly scraped fixtures as rows 295 and 296.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F..AC (i.e. everything except A, C, D, E) as 1-based column indices.
$cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

function Swap-Rows($r1, $r2) {
    foreach ($c in $cols) {
        $v1 = $ws.Cells.Item($r1, $c).Value()
        $v2 = $ws.Cells.Item($r2, $c).Value()
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

Swap-Rows 55 56
Swap-Rows 131 132
Swap-Rows 221 222
Swap-Rows 229 231
Swap-Rows 245 246
Swap-Rows 263 264
Swap-Rows 271 272

# ---------------------------------------------------------------------
# Append two new fixture rows at the bottom of the sheet (295 and 296).
# Neither match has been played yet, so FTHG/FTAG/FTR (H, I, J) and the
# closing handicap/over-under profit columns (AB, AC) are left blank.
# ---------------------------------------------------------------------

function Set-Row($r, $vals) {
    foreach ($key in $vals.Keys) {
        $ws.Range("$key$r").Value = $vals[$key]
    }
}

$row295 = @{
    "A" = 293
    "B" = 7811117
    "C" = "Mexico Liga MX Femenil"
    "D" = "Mexico Liga MX Femenil"
    "E" = 45401
    "F" = "Monterrey Women"
    "G" = "Toluca Women"
    "K" = 1.4
    "L" = 4.6
    "M" = 5.5
    "N" = 1.111
    "O" = 7.5
    "P" = 19
    "Q" = -2.25
    "R" = 2
    "S" = 1.8
    "T" = 3.5
    "U" = 1.9
    "V" = 1.9
    "W" = 0
    "X" = 0
    "Y" = 0
    "Z" = 0
    "AA" = 0
}

$row296 = @{
    "A" = 294
    "B" = 7645828
    "C" = "Mexico Liga MX Femenil"
    "D" = "Mexico Liga MX Femenil"
    "E" = 45401.83333333334
    "F" = "Queretaro Women"
    "G" = "Leon Women"
    "K" = 2.35
    "L" = 3.6
    "M" = 2.45
    "N" = 2.4
    "O" = 3.6
    "P" = 2.4
    "Q" = 0
    "R" = 1.925
    "S" = 1.875
    "T" = 2.75
    "U" = 1.95
    "V" = 1.85
    "W" = 0
    "X" = 0
    "Y" = 0
    "Z" = 0
    "AA" = 0
}

Set-Row 295 $row295
Set-Row 296 $row296

# Match the bold/bordered/centered style used on column A (id) and the
# custom date/time number format used on column E (Date) elsewhere in
# the sheet, by cloning the formatting from the row directly above.
$ws.Range("A294").Copy()
$ws.Range("A295").PasteSpecial(-4122)
$ws.Range("A294").Copy()
$ws.Range("A296").PasteSpecial(-4122)

$ws.Range("E294").Copy()
$ws.Range("E295").PasteSpecial(-4122)
$ws.Range("E294").Copy()
$ws.Range("E296").PasteSpecial(-4122)

Write-Host "Edit applied"
